$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lama1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06628666666666667
$ws.Range("H2").Value = 0.19886
$ws.Range("I2").Value = 0.09845363529874238
$ws.Range("J2").Value = 0.09845363529874238
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 10.57163870670222
$ws.Range("R2").Value = 95.14474836032001
$ws.Range("S2").Value = 0.02937152617411446
$ws.Range("T2").Value = 0.02937152617411446

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lama1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06628666666666667
$ws.Range("H3").Value = 0.19886
$ws.Range("I3").Value = 0.09845363529874238
$ws.Range("J3").Value = 0.09845363529874238
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 11.43832757114
$ws.Range("R3").Value = 102.94494814026
$ws.Range("S3").Value = 0.0317794759133076
$ws.Range("T3").Value = 0.0317794759133076

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lama1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06628666666666667
$ws.Range("H4").Value = 0.19886
$ws.Range("I4").Value = 0.09845363529874238
$ws.Range("J4").Value = 0.09845363529874238
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 4.930913049624444
$ws.Range("R4").Value = 44.37821744662
$ws.Range("S4").Value = 0.0136997154100157
$ws.Range("T4").Value = 0.0136997154100157

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lama1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.06628666666666667
$ws.Range("H5").Value = 0.19886
$ws.Range("I5").Value = 0.09845363529874238
$ws.Range("J5").Value = 0.09845363529874238
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 3.872110068775556
$ws.Range("R5").Value = 34.84899061898
$ws.Range("S5").Value = 0.01075800879971341
$ws.Range("T5").Value = 0.01075800879971341

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Lama1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.06628666666666667
$ws.Range("H6").Value = 0.19886
$ws.Range("I6").Value = 0.09845363529874238
$ws.Range("J6").Value = 0.09845363529874238
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 4.623244171253334
$ws.Range("R6").Value = 41.60919754128
$ws.Range("S6").Value = 0.0128449090015912
$ws.Range("T6").Value = 0.01284490900159121

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lama1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.585754
$ws.Range("H7").Value = 1.757262
$ws.Range("I7").Value = 0.8700031784790234
$ws.Range("J7").Value = 0.8700031784790236
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 93.41817850254934
$ws.Range("R7").Value = 840.7636065229441
$ws.Range("S7").Value = 0.2595467506174028
$ws.Range("T7").Value = 0.2595467506174028

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lama1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.585754
$ws.Range("H8").Value = 1.757262
$ws.Range("I8").Value = 0.8700031784790234
$ws.Range("J8").Value = 0.8700031784790236
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 101.076829851738
$ws.Range("R8").Value = 909.691468665642
$ws.Range("S8").Value = 0.2808250296810355
$ws.Range("T8").Value = 0.2808250296810356

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lama1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.585754
$ws.Range("H9").Value = 1.757262
$ws.Range("I9").Value = 0.8700031784790234
$ws.Range("J9").Value = 0.8700031784790236
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 43.57289614507266
$ws.Range("R9").Value = 392.156065305654
$ws.Range("S9").Value = 0.1210599884382732
$ws.Range("T9").Value = 0.1210599884382732

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Lama1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.585754
$ws.Range("H10").Value = 1.757262
$ws.Range("I10").Value = 0.8700031784790234
$ws.Range("J10").Value = 0.8700031784790236
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 34.21659400420733
$ws.Range("R10").Value = 307.949346037866
$ws.Range("S10").Value = 0.09506507120286622
$ws.Range("T10").Value = 0.09506507120286625

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Lama1"
$ws.Range("C11").Value = "Itgb1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.585754
$ws.Range("H11").Value = 1.757262
$ws.Range("I11").Value = 0.8700031784790234
$ws.Range("J11").Value = 0.8700031784790236
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 40.854125006864
$ws.Range("R11").Value = 367.687125061776
$ws.Range("S11").Value = 0.1135063385394456
$ws.Range("T11").Value = 0.1135063385394457

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Lama1"
$ws.Range("C12").Value = "Itgb1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.01811333333333333
$ws.Range("H12").Value = 0.05434
$ws.Range("I12").Value = 0.02690320095611817
$ws.Range("J12").Value = 0.02690320095611818
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 2.888780284231111
$ws.Range("R12").Value = 25.99902255808
$ws.Range("S12").Value = 0.00802599181485155
$ws.Range("T12").Value = 0.008025991814851552

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Lama1"
$ws.Range("C13").Value = "Itgb1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.01811333333333333
$ws.Range("H13").Value = 0.05434
$ws.Range("I13").Value = 0.02690320095611817
$ws.Range("J13").Value = 0.02690320095611818
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 3.12560957566
$ws.Range("R13").Value = 28.13048618094
$ws.Range("S13").Value = 0.008683982304782935
$ws.Range("T13").Value = 0.008683982304782937

# Row 14
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Lama1"
$ws.Range("C14").Value = "Itgb1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.01811333333333333
$ws.Range("H14").Value = 0.05434
$ws.Range("I14").Value = 0.02690320095611817
$ws.Range("J14").Value = 0.02690320095611818
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 1.347409308642222
$ws.Range("R14").Value = 12.12668377778
$ws.Range("S14").Value = 0.0037435509171289
$ws.Range("T14").Value = 0.003743550917128901

# Row 15
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Lama1"
$ws.Range("C15").Value = "Itgb1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.01811333333333333
$ws.Range("H15").Value = 0.05434
$ws.Range("I15").Value = 0.02690320095611817
$ws.Range("J15").Value = 0.02690320095611818
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 1.058083380957778
$ws.Range("R15").Value = 9.52275042862
$ws.Range("S15").Value = 0.002939707322621072
$ws.Range("T15").Value = 0.002939707322621073

# Row 16
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Lama1"
$ws.Range("C16").Value = "Itgb1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.01811333333333333
$ws.Range("H16").Value = 0.05434
$ws.Range("I16").Value = 0.02690320095611817
$ws.Range("J16").Value = 0.02690320095611818
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 1.263336459146667
$ws.Range("R16").Value = 11.37002813232
$ws.Range("S16").Value = 0.003509968596733712
$ws.Range("T16").Value = 0.003509968596733713

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Lama1"
$ws.Range("C17").Value = "Itgb1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.003124
$ws.Range("H17").Value = 0.009372
$ws.Range("I17").Value = 0.004639985266115928
$ws.Range("J17").Value = 0.004639985266115929
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 0.4982268830293334
$ws.Range("R17").Value = 4.484041947264
$ws.Range("S17").Value = 0.001384239883856988
$ws.Range("T17").Value = 0.001384239883856988

# Row 18
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Lama1"
$ws.Range("C18").Value = "Itgb1"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.003124
$ws.Range("H18").Value = 0.009372
$ws.Range("I18").Value = 0.004639985266115928
$ws.Range("J18").Value = 0.004639985266115929
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 0.5390727446279999
$ws.Range("R18").Value = 4.851654701652
$ws.Range("S18").Value = 0.001497723263901834
$ws.Range("T18").Value = 0.001497723263901834

# Row 19
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Lama1"
$ws.Range("C19").Value = "Itgb1"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.003124
$ws.Range("H19").Value = 0.009372
$ws.Range("I19").Value = 0.004639985266115928
$ws.Range("J19").Value = 0.004639985266115929
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 0.2323871925026667
$ws.Range("R19").Value = 2.091484732524
$ws.Range("S19").Value = 0.0006456488626303286
$ws.Range("T19").Value = 0.0006456488626303288

# Row 20
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Lama1"
$ws.Range("C20").Value = "Itgb1"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.003124
$ws.Range("H20").Value = 0.009372
$ws.Range("I20").Value = 0.004639985266115928
$ws.Range("J20").Value = 0.004639985266115929
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 0.1824872551773333
$ws.Range("R20").Value = 1.642385296596
$ws.Range("S20").Value = 0.0005070102507840392
$ws.Range("T20").Value = 0.0005070102507840394

# Row 21
$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Lama1"
$ws.Range("C21").Value = "Itgb1"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.003124
$ws.Range("H21").Value = 0.009372
$ws.Range("I21").Value = 0.004639985266115928
$ws.Range("J21").Value = 0.004639985266115929
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 0.217887178784
$ws.Range("R21").Value = 1.960984609056
$ws.Range("S21").Value = 0.0006053630049427374
$ws.Range("T21").Value = 0.0006053630049427377
